$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.501.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.94%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.164.36'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.16%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.16%  '

# Row 6
$ws.Range("E6").Value = '  +0.97%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '64.07'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.91%  '

# Row 8
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.397'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0858'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.05%  '

# Row 11
$ws.Range("E11").Value = '  -0.31%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.27'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.73%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.484.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.13%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.26%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.815'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.17%  '

# Row 16
$ws.Range("E16").Value = '  +0.48%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.162.06'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.29%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.456.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.84%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.06%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.15'
$ws.Range("D20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = '  +1.53%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '230.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.95%  '

# Row 24
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.85%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.66%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.69'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.26%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.49'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.39%  '

# Row 28
$ws.Range("E28").Value = '  +2.09%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.00%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.43'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.63%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.83%  '

# Row 32
$ws.Range("E32").Value = '  +1.43%  '

# Row 33
$ws.Range("E33").Value = '  +1.30%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.41%  '

# Row 35
$ws.Range("E35").Value = '  -0.95%  '

# Row 36
$ws.Range("E36").Value = '  -0.48%  '

# Row 37
$ws.Range("E37").Value = '  +1.25%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.57'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.21%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.999'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.18%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '103.56'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.19%  '

# Row 41
$ws.Range("E41").Value = '  +0.75%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.06%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.530.62'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.25%  '

# Row 44
$ws.Range("E44").Value = '  +4.22%  '

# Row 45
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.34'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.63%  '

# Row 46
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0934'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.61%  '

# Row 47
$ws.Range("E47").Value = '  +0.57%  '

# Row 48
$ws.Range("E48").Value = '  +5.25%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.367.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.27%  '

# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.96'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.39%  '
